# Generate Report for Handoff
# The c0354d71-4316-4e93-9db7-e901cb37535e.md file's handback is now stale
# (outdated vs. latest source), so it moves from "Handed back: in sync with
# en-US" to "Ready for handoff", with updated handoff timestamps and a new
# error detail message, across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$statusReadyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1693edd0fa96772f57a670ff2ac3b73b1336861d/e2e/c0354d71-4316-4e93-9db7-e901cb37535e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/789cb6e046d3040a65814e08255af60d2e7836be/e2e/c0354d71-4316-4e93-9db7-e901cb37535e.md."

# ----- Overview sheet -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReadyForHandoff
$wsOverview.Range("F3").Value = $statusReadyForHandoff
$wsOverview.Range("G3").Value = "2016-08-24 18:57:29"

# ----- zh-cn sheet -----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReadyForHandoff
$wsZhCn.Range("H3").Value = "2016-08-24 18:57:24"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ----- de-de sheet -----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReadyForHandoff
$wsDeDe.Range("H3").Value = "2016-08-24 18:57:29"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
